$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove the last data row (old row 60 -> "Winnona Adult"),
# since the data shifted up by one row and the final row dropped off.
$ws.Rows(60).Delete()

# Step 2: update the cells whose values changed.
$ws.Range("H2").Value = 0.5275666176396232
$ws.Range("H3").Value = 1.111453171863835
$ws.Range("H4").Value = 0.4893562365297635
$ws.Range("H5").Value = 1.991946537346263
$ws.Range("H6").Value = 1.39505162991666
$ws.Range("H7").Value = 1.409308091756941
$ws.Range("H8").Value = 1.137786789283417
$ws.Range("H9").Value = 0.9986028064211235
$ws.Range("H10").Value = 1.076576664650283
$ws.Range("H11").Value = 0.7514764911590116
$ws.Range("H12").Value = 1.070214770847774
$ws.Range("H13").Value = 0.6970488565016943
$ws.Range("H14").Value = 1.391632043822133
$ws.Range("H15").Value = 0.9969763735310998
$ws.Range("H16").Value = 0.8203147244775868
$ws.Range("H17").Value = 0.753798502680639
$ws.Range("H18").Value = 1.591646543446408
$ws.Range("H19").Value = 2.069665260314198
$ws.Range("H20").Value = 0.7649637553464128
$ws.Range("H21").Value = 1.615221989336772
$ws.Range("H22").Value = 2.100321112617971
$ws.Range("H23").Value = 1.201519058846645
$ws.Range("H24").Value = 1.329752906257196
$ws.Range("H25").Value = 0.4893259532382404
$ws.Range("H26").Value = 1.091034515845942
$ws.Range("H27").Value = 0.4538852514670935
$ws.Range("H28").Value = 1.418704576432297
$ws.Range("H29").Value = 1.758336896544703
$ws.Range("H30").Value = 1.011781718944878
$ws.Range("H35").Value = 0.604890884791715
$ws.Range("H36").Value = 1.290839341586115
$ws.Range("H37").Value = 1.678516724035213
$ws.Range("H38").Value = 0.8174928716171911
$ws.Range("E39").Value = "Adult"
$ws.Range("H39").Value = 1.305641572198231
$ws.Range("B40").Value = 26
$ws.Range("C40").Value = "North Memorial Adult"
$ws.Range("D40").Value = "North Memorial Health Hospital"
$ws.Range("F40").Value = 0.2107983917288914
$ws.Range("G40").Value = 333.1226158038147
$ws.Range("H40").Value = 0.7897727313129739
$ws.Range("B41").Value = 10
$ws.Range("C41").Value = "Owatonna Adult"
$ws.Range("D41").Value = "Owatonna Hospital"
$ws.Range("F41").Value = 0.1728891441700172
$ws.Range("G41").Value = 217.6744186046512
$ws.Range("H41").Value = 1.208645277805187
$ws.Range("B42").Value = 71
$ws.Range("C42").Value = "PrairieCare Pediatric/Adolescent"
$ws.Range("D42").Value = "PrairieCare"
$ws.Range("E42").Value = "Adolescent"
$ws.Range("F42").Value = 3.429611650485437
$ws.Range("G42").Value = 184.5859872611465
$ws.Range("H42").Value = 0.6750192592315432
$ws.Range("E43").Value = "Child"
$ws.Range("H43").Value = 0.6261292379729364
$ws.Range("B44").Value = 14
$ws.Range("C44").Value = "Regina Geriatric"
$ws.Range("D44").Value = "Regina Hospital"
$ws.Range("E44").Value = "Geriatric"
$ws.Range("F44").Value = 0.02067777139574957
$ws.Range("G44").Value = 572
$ws.Range("H44").Value = 0.5980860922294489
$ws.Range("B45").Value = 100
$ws.Range("C45").Value = "Regions Adult"
$ws.Range("D45").Value = "Regions Hospital"
$ws.Range("E45").Value = "Adult"
$ws.Range("F45").Value = 1.136128661688685
$ws.Range("G45").Value = 267.360970677452
$ws.Range("H45").Value = 0.9840297836997983
$ws.Range("B46").Value = 12
$ws.Range("C46").Value = "Ridgeview Geriatric"
$ws.Range("D46").Value = "Ridgeview Medical Center"
$ws.Range("E46").Value = "Geriatric"
$ws.Range("F46").Value = 0.007466973004020678
$ws.Range("G46").Value = 164.3076923076923
$ws.Range("H46").Value = 2.082101208716377
$ws.Range("B47").Value = 37
$ws.Range("C47").Value = "St Jospephs Adult/Geriatric"
$ws.Range("D47").Value = "Saint Joseph's Hospital"
$ws.Range("E47").Value = "Adult"
$ws.Range("F47").Value = 0.7375071797817346
$ws.Range("G47").Value = 262.6542056074766
$ws.Range("H47").Value = 1.001663603813287
$ws.Range("E48").Value = "Geriatric"
$ws.Range("H48").Value = 1.302492925875719
$ws.Range("B49").Value = 16
$ws.Range("C49").Value = "Sanford Behavioral Adolescent/Adult/Geriatric"
$ws.Range("D49").Value = "Sanford Behavioral Health Center"
$ws.Range("E49").Value = "Adolescent"
$ws.Range("F49").Value = 0.08836123005060335
$ws.Range("G49").Value = 145.057268722467
$ws.Range("H49").Value = 0.8589648590718558
$ws.Range("E50").Value = "Adult"
$ws.Range("H50").Value = 1.813705445184304
$ws.Range("E51").Value = "Geriatric"
$ws.Range("H51").Value = 2.358415043714789
$ws.Range("B52").Value = 28
$ws.Range("C52").Value = "St. Cloud Adolescent/Adult"
$ws.Range("D52").Value = "St. Cloud Hospital"
$ws.Range("E52").Value = "Adolescent"
$ws.Range("F52").Value = 0.3386531724406384
$ws.Range("G52").Value = 180
$ws.Range("H52").Value = 0.692217202141901
$ws.Range("E53").Value = "Adult"
$ws.Range("H53").Value = 1.461617545252784
$ws.Range("B54").Value = 22
$ws.Range("C54").Value = "Essentia St. Josephs Adult"
$ws.Range("D54").Value = "St. Joseph's Medical Center"
$ws.Range("F54").Value = 0.2056289488799541
$ws.Range("G54").Value = 124.6256983240224
$ws.Range("H54").Value = 2.111050623455474
$ws.Range("C55").Value = "St Lukes Adult"
$ws.Range("D55").Value = "St. Luke's Hospital"
$ws.Range("F55").Value = 0.2952326249282022
$ws.Range("G55").Value = 268.0622568093385
$ws.Range("H55").Value = 0.9814554323200634
$ws.Range("B56").Value = 16
$ws.Range("C56").Value = "United Adolescent"
$ws.Range("D56").Value = "United Hospital"
$ws.Range("E56").Value = "Adolescent"
$ws.Range("F56").Value = 0.299757281553398
$ws.Range("G56").Value = 258.6558704453441
$ws.Range("H56").Value = 0.4817176434890577
$ws.Range("B57").Value = 42
$ws.Range("C57").Value = "United Adult/Geriatric"
$ws.Range("E57").Value = "Adult"
$ws.Range("F57").Value = 0.3739230327398047
$ws.Range("G57").Value = 240.073732718894
$ws.Range("H57").Value = 1.095876484136474
$ws.Range("E58").Value = "Geriatric"
$ws.Range("H58").Value = 1.425000731570334
$ws.Range("B59").Value = 10
$ws.Range("C59").Value = "Winnona Adult"
$ws.Range("D59").Value = "Winona Health Services"
$ws.Range("E59").Value = "Adult"
$ws.Range("F59").Value = 0.01033888569787478
$ws.Range("G59").Value = 169.3333333333333
$ws.Range("H59").Value = 1.553687941804141
